$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "barren soil/land" terminology to simpler "ground/land" wording
$ws.Range("A12").Value = "On dry ground"
$ws.Range("A14").Value = "On wet ground"
$ws.Range("B14").Value = "Towards wet land"
$ws.Range("B12").Value = "Towards dry land"

$ws.Range("C13").Value = "On dry ground"
$ws.Range("D13").Value = "On dry ground"

$ws.Range("C15").Value = "On wet ground"
$ws.Range("D15").Value = "On wet ground;On a frozen snowland;In a lightly forested area;On the side of a stream;On the shore of a lake;On frozen land;On the shore of a river;On the shores of an ocean;Swamp"

# Add a second required-previous-place-name entry for the Swamp row (2 destinations max without a map)
$ws.Range("C19").Value = "On wet ground"

# Update the selected cell in the sheet view
$ws.Range("C19").Select()
